# Sprint Backlog / Burndown workbook update
# - Fill in the Week 2 (and Week 1/3/4) actuals for each task row (D3:G9)
# - Recalculated totals row (C27:G27) / burndown chart pick this up automatically
#   via the existing shared formula in row 27.
# - Update the saved cursor/selection to match the author's last position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task rows: Week1 (D), Week2 (E), Week3 (F), Week4 (G) -----------------
# Row 3: Implement UI (fxml and codebehind)
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

# Row 4: Implement storing location information
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Row 5: Design world
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Row 6: Implement loading world information
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# Row 7: Implement move action & player
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

# Row 8: UI binding to take action & update display
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# Row 9: Prepare and make submission
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.5
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# --- View state: scroll position + active selection ------------------------
$win = $excel.ActiveWindow
try {
    $win.ScrollColumn = 3
    $win.ScrollRow = 4
} catch {
}

$ws.Range("E11").Select()
